$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = $false
$ws.Range("E2").Value = 77.349997999999999
$ws.Range("F2").Value = -1.0869590792838919
$ws.Range("G2").Value = $false

# New row 3
$ws.Range("C3").Value = 9891.2999999999993

# Column width adjustments (closest achievable widths given engine's width quantization)
$ws.Range("E1").ColumnWidth = 9.0
$ws.Range("F1").ColumnWidth = 11.666666666666666
